$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header text: P1 (syst3_c -> syst3_u) then O1 (syst2_c -> syst2_u)
# Order matters for how new shared-string entries are appended.
$ws.Range("P1").Value = "syst3_u"
$ws.Range("O1").Value = "syst2_u"

# Move active cell selection to O16 (matches final sheetView selection)
$ws.Range("O16").Select()
